$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value parses as a plain number must be forced back to Text,
# matching the source data model where the Price column is stored as inline strings
# (e.g. "67.697.48"), not numeric cells. We briefly set NumberFormat to Text ("@"),
# write the value, then restore the cell style so no stray formatting is left behind.

$ws.Range("D2").Value = "67.697.48"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.491.77"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("E9").Value = "  +5.20%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").Value = "2.922.34"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "67.602.34"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "2.488.44"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("E21").Value = "  +2.26%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.27%  "
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.75%  "
$ws.Range("D27").Value = "2.620.54"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").Value = [string]::Concat("0.0", [char]0x2083, "0910")
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "510.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  +4.13%  "
$ws.Range("E32").Value = "  +3.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "161.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.06%  "
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.60%  "
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "144.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = [string]::Concat("0.0", [char]0x2086, "0258")
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.514"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.585"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.24%  "
